# Append: 2025-12-22 01:28 JST
#
# A fresh scrape added one new listing ("WebRTCを用いたビデオ通話サイトの制作依頼")
# above the previously-first "初回 ECサイト要件定義..." entry, pushing every
# row from the old row 8 down by one, and refreshed the "取得日時" capture
# timestamp (column A) on every listing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-22 01:28:15"

# Make room for the new listing by inserting a blank row at row 8; this
# shifts the old rows 8-11 down to rows 9-12.
$ws.Cells.Item(8, 1).EntireRow.Insert()

# --- Fill every data row's "取得日時" column with the refreshed timestamp ---
$ws.Range("A2:A12").Value = $newTimestamp

# --- Populate the newly inserted row 8 with the new listing's data ---
$ws.Range("B8").Value = "【急募】WebRTCを用いたビデオ通話サイトの制作依頼"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5458447"
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = "◇サイト"

# --- Rebuild the F-column hyperlinks so they line up with their row again ---
# (Row insertion does not re-anchor existing hyperlink objects to the
# content that slid underneath them, so drop them all and re-add in the
# current top-to-bottom order; each target now matches its row's URL text.)
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5458419",
    "https://www.lancers.jp/work/detail/5458381",
    "https://www.lancers.jp/work/detail/5458190",
    "https://www.lancers.jp/work/detail/5458166",
    "https://www.lancers.jp/work/detail/5458299",
    "https://www.lancers.jp/work/detail/5431107",
    "https://www.lancers.jp/work/detail/5458447",
    "https://www.lancers.jp/work/detail/5425629",
    "https://www.lancers.jp/work/detail/5458330",
    "https://www.lancers.jp/work/detail/5458234",
    "https://www.lancers.jp/work/detail/5458288"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}
